# Refresh Leve profit-tracking figures (currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets with the latest
# Universalis market-board snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

# ==================== Sheet: ALC ====================
$ws = $wb.Worksheets.Item(1)
# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 1224.8125
$ws.Range("I62").Value = 1211.875
$ws.Range("J62").Value = 1237.75
$ws.Range("K62").Value = 1211.875
$ws.Range("L62").Value = 1237.75
$ws.Range("M62").Value = -587.875
$ws.Range("N62").Value = -2485.75

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 1224.8125
$ws.Range("I65").Value = 1211.875
$ws.Range("J65").Value = 1237.75
$ws.Range("K65").Value = 6059.375
$ws.Range("L65").Value = 6188.75
$ws.Range("M65").Value = -2939.375
$ws.Range("N65").Value = -12428.75

# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 2759.4
$ws.Range("I70").Value = 3084.5715
$ws.Range("J70").Value = 2474.875
$ws.Range("K70").Value = 9253.7145
$ws.Range("L70").Value = 7424.625
$ws.Range("M70").Value = -8983.7145
$ws.Range("N70").Value = -7964.625

# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 2759.4
$ws.Range("I73").Value = 3084.5715
$ws.Range("J73").Value = 2474.875
$ws.Range("K73").Value = 9253.7145
$ws.Range("L73").Value = 7424.625
$ws.Range("M73").Value = -8317.7145
$ws.Range("N73").Value = -9296.625

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 919.2059
$ws.Range("I137").Value = 796.8253999999999
$ws.Range("K137").Value = 2390.4762
$ws.Range("M137").Value = 159.5237999999999

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 4681.7646
$ws.Range("I138").Value = 1050.16
$ws.Range("J138").Value = 8173.6924
$ws.Range("K138").Value = 3150.48
$ws.Range("L138").Value = 24521.0772
$ws.Range("M138").Value = 1989.52
$ws.Range("N138").Value = -34801.0772

# ==================== Sheet: ARM ====================
$ws = $wb.Worksheets.Item(2)
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1224.3334
$ws.Range("I74").Value = 1099.0741
$ws.Range("J74").Value = 1600.1111
$ws.Range("K74").Value = 1099.0741
$ws.Range("L74").Value = 1600.1111
$ws.Range("M74").Value = -225.0741
$ws.Range("N74").Value = -3348.1111

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1224.3334
$ws.Range("I77").Value = 1099.0741
$ws.Range("J77").Value = 1600.1111
$ws.Range("K77").Value = 5495.3705
$ws.Range("L77").Value = 8000.5555
$ws.Range("M77").Value = -1127.3705
$ws.Range("N77").Value = -16736.5555

# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 2847.5
$ws.Range("I88").Value = 2793.25
$ws.Range("J88").Value = 2901.75
$ws.Range("K88").Value = 2793.25
$ws.Range("L88").Value = 2901.75
$ws.Range("M88").Value = -2387.25
$ws.Range("N88").Value = -3713.75

# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 2847.5
$ws.Range("I91").Value = 2793.25
$ws.Range("J91").Value = 2901.75
$ws.Range("K91").Value = 2793.25
$ws.Range("L91").Value = 2901.75
$ws.Range("M91").Value = -1389.25
$ws.Range("N91").Value = -5709.75

# Row 123 (Leve Item ID 34107)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2321.2131
$ws.Range("I132").Value = 1519.775
$ws.Range("J132").Value = 3847.762
$ws.Range("K132").Value = 4559.325000000001
$ws.Range("L132").Value = 11543.286
$ws.Range("M132").Value = -2029.325000000001
$ws.Range("N132").Value = -16603.286

# Row 140 (Leve Item ID 42496)
$ws.Range("H140").Value = 46466.668
$ws.Range("J140").Value = 46466.668
$ws.Range("L140").Value = 46466.668
$ws.Range("N140").Value = -56826.668

# ==================== Sheet: BSM ====================
$ws = $wb.Worksheets.Item(3)
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 1018.5
$ws.Range("I94").Value = 665.9286
$ws.Range("J94").Value = 2252.5
$ws.Range("K94").Value = 665.9286
$ws.Range("L94").Value = 2252.5
$ws.Range("M94").Value = -214.9286
$ws.Range("N94").Value = -3154.5

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 4201.5713
$ws.Range("I134").Value = 4908.4
$ws.Range("J134").Value = 2434.5
$ws.Range("K134").Value = 14725.2
$ws.Range("L134").Value = 7303.5
$ws.Range("M134").Value = -12190.2
$ws.Range("N134").Value = -12373.5

# ==================== Sheet: CRP ====================
$ws = $wb.Worksheets.Item(4)
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1824.8334
$ws.Range("I16").Value = 1642.5714
$ws.Range("J16").Value = 2080
$ws.Range("K16").Value = 1642.5714
$ws.Range("L16").Value = 2080
$ws.Range("M16").Value = -1355.5714
$ws.Range("N16").Value = -2654

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2624.93
$ws.Range("I31").Value = 1246.5588
$ws.Range("K31").Value = 1246.5588
$ws.Range("M31").Value = -951.5588

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2624.93
$ws.Range("I34").Value = 1246.5588
$ws.Range("K34").Value = 1246.5588
$ws.Range("M34").Value = -1044.5588

# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1824.8334
$ws.Range("I113").Value = 1642.5714
$ws.Range("J113").Value = 2080
$ws.Range("K113").Value = 1642.5714
$ws.Range("L113").Value = 2080
$ws.Range("M113").Value = 527.4286
$ws.Range("N113").Value = -6420

# Row 135 (Leve Item ID 42008)
$ws.Range("H135").Value = 32884
$ws.Range("J135").Value = 32884
$ws.Range("L135").Value = 32884
$ws.Range("N135").Value = -43024

# ==================== Sheet: CUL ====================
$ws = $wb.Worksheets.Item(5)
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1516082.5
$ws.Range("I131").Value = 5882952
$ws.Range("J131").Value = 1046.0408
$ws.Range("K131").Value = 17648856
$ws.Range("L131").Value = 3138.1224
$ws.Range("M131").Value = -17643816
$ws.Range("N131").Value = -13218.1224

# ==================== Sheet: GSM ====================
$ws = $wb.Worksheets.Item(6)
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 891.5
$ws.Range("I97").Value = 891.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 891.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -395.5
$ws.Range("N97").ClearContents()

# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 22887.727
$ws.Range("J123").Value = 22887.727
$ws.Range("L123").Value = 22887.727
$ws.Range("N123").Value = -27787.727

# Row 128 (Leve Item ID 34544)
$ws.Range("H128").Value = 52780
$ws.Range("J128").Value = 52780
$ws.Range("L128").Value = 52780
$ws.Range("N128").Value = -62740

# ==================== Sheet: LTW ====================
$ws = $wb.Worksheets.Item(7)
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 71431090
$ws.Range("I68").Value = 2267.5557
$ws.Range("J68").Value = 200002960
$ws.Range("K68").Value = 2267.5557
$ws.Range("L68").Value = 200002960
$ws.Range("M68").Value = -1518.5557
$ws.Range("N68").Value = -200004458

# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 71431090
$ws.Range("I71").Value = 2267.5557
$ws.Range("J71").Value = 200002960
$ws.Range("K71").Value = 11337.7785
$ws.Range("L71").Value = 1000014800
$ws.Range("M71").Value = -7593.7785
$ws.Range("N71").Value = -1000022288

# ==================== Sheet: WVR ====================
$ws = $wb.Worksheets.Item(8)
# Row 4 (Leve Item ID 2996)
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 2489.077
$ws.Range("I96").Value = 2101.625
$ws.Range("K96").Value = 2101.625
$ws.Range("M96").Value = -728.625
